$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clear out the old scratch calc block (J11:K14) that is gone in the
#    new version of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("J11:K14").ClearContents()

# ---------------------------------------------------------------------------
# 2. Currency number format used throughout the Cost / Cost Individual
#    columns (builtin format id 8: "$"#,##0.00_);[Red]("$"#,##0.00))
# ---------------------------------------------------------------------------
$currencyFmt = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

# ---------------------------------------------------------------------------
# 3. Row 3 - 555 timer (hyperlink added first so the workbook-level
#    "Hyperlink" cell style picks up the theme-based link color; B2 below
#    then reuses that same style without keeping a live hyperlink).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "555 timer"
$ws.Range("B3").Value = "https://www.digikey.com/en/products/detail/texas-instruments/NE555P/277057"
$ws.Range("C3").Value = 0.48
$ws.Range("C3").NumberFormat = $currencyFmt
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.48
$ws.Range("E3").NumberFormat = $currencyFmt
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.digikey.com/en/products/detail/texas-instruments/NE555P/277057")

# ---------------------------------------------------------------------------
# 4. Row 2 - Op Amp (link-styled text, but not an actual clickable hyperlink)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Op Amp"
$ws.Range("B2").Value = "https://www.digikey.com/en/products/detail/texas-instruments/LMC660CN-NOPB/32519"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.digikey.com/en/products/detail/texas-instruments/LMC660CN-NOPB/32519") | Out-Null
$ws.Hyperlinks.Item(2).Delete()
$ws.Range("C2").Value = 3.06
$ws.Range("C2").NumberFormat = $currencyFmt
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 6.12
$ws.Range("E2").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# 5. Row 4 - Diodes
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Diodes"
$ws.Range("B4").Value = "https://www.digikey.com/en/products/detail/smc-diode-solutions/SF26G/6022671"
$ws.Range("C4").Value = 0.38
$ws.Range("C4").NumberFormat = $currencyFmt
$ws.Range("D4").Value = 4
$ws.Range("E4").Formula = "=(C4)*4"
$ws.Range("E4").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# 6. Row 5 - Transformer
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Transformer"
$ws.Range("B5").Value = "https://www.digikey.com/en/products/detail/hammond-manufacturing/166LA12/2182731"
$ws.Range("C5").Value = 17.72
$ws.Range("C5").NumberFormat = $currencyFmt
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 17.72
$ws.Range("E5").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# 7. Row 6 - Capacitor
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Capacitor"
$ws.Range("B6").Value = "https://www.digikey.com/en/products/detail/cornell-dubilier-illinois-capacitor/109LBB016M2BC/5410933"
$ws.Range("C6").Value = 3.47
$ws.Range("C6").NumberFormat = $currencyFmt
$ws.Range("D6").Value = 2
$ws.Range("E6").Formula = "=C6*2"
$ws.Range("E6").NumberFormat = $currencyFmt
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.digikey.com/en/products/detail/cornell-dubilier-illinois-capacitor/109LBB016M2BC/5410933")

# ---------------------------------------------------------------------------
# 8. Row 7 - other capacitor
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "other capacitor"
$ws.Range("B7").Value = "https://www.digikey.com/en/products/detail/cornell-dubilier-illinois-capacitor/105CKH050M/5410526"
$ws.Range("C7").Value = 0.27
$ws.Range("C7").NumberFormat = $currencyFmt
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.27
$ws.Range("E7").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# 9. Row 8 - one last capacitor
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "one last capacitor"
$ws.Range("B8").Value = "https://www.digikey.com/en/products/detail/nichicon/URZ2AR22MDD1TD/4320687"
$ws.Range("C8").Value = 0.33
$ws.Range("C8").NumberFormat = $currencyFmt
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.33
$ws.Range("E8").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# 10. Row 9 - resistors
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "resistors"

# ---------------------------------------------------------------------------
# 11. Selection / active cell, matching the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("C14").Select()
